$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared string "Time 24 small" -> "Time 24 hrs" (cell A3)
$ws.Range("A3").Value = "Time 24 hrs"

# Update B3 value (combined amazon size fractions)
$ws.Range("B3").Value = 0.065740715048571

# Bold the header cell A1 (creates a new cellXf with applyFont=true)
$ws.Range("A1").Font.Bold = $true

# Update active cell selection to B7
$ws.Range("B7").Select()
